$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 510266.4
$ws.Range("I9").Value = 14551.714
$ws.Range("J9").Value = 1666934
$ws.Range("K9").Value = 14551.714
$ws.Range("L9").Value = 1666934
$ws.Range("M9").Value = -14382.714
$ws.Range("N9").Value = -1667272
$ws.Range("H28").Value = 228
$ws.Range("I28").Value = 256.0625
$ws.Range("K28").Value = 256.0625
$ws.Range("M28").Value = 228.9375
$ws.Range("H40").Value = 3812
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H41").Value = 239.04546
$ws.Range("I41").Value = 208.05
$ws.Range("K41").Value = 208.05
$ws.Range("M41").Value = 231.95
$ws.Range("H106").Value = 5438
$ws.Range("I106").Value = 5616.6665
$ws.Range("K106").Value = 5616.6665
$ws.Range("M106").Value = -4985.6665
$ws.Range("H135").Value = 1450.0834
$ws.Range("I135").Value = 1450.0834
$ws.Range("K135").Value = 13050.7506
$ws.Range("M135").Value = -10515.7506
$ws.Range("H137").Value = 4209.283
$ws.Range("I137").Value = 2331.6
$ws.Range("K137").Value = 6994.799999999999
$ws.Range("M137").Value = -4444.799999999999
$ws.Range("H138").Value = 2733.57
$ws.Range("I138").Value = 2267.4
$ws.Range("J138").Value = 2758.1052
$ws.Range("K138").Value = 6802.200000000001
$ws.Range("L138").Value = 8274.3156
$ws.Range("M138").Value = -1662.200000000001
$ws.Range("N138").Value = -18554.3156

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3604.7778
$ws.Range("J61").Value = 4382.0557
$ws.Range("L61").Value = 4382.0557
$ws.Range("N61").Value = -4806.0557
$ws.Range("H110").Value = 3952.95
$ws.Range("I110").Value = 3877.7334
$ws.Range("J110").Value = 4178.6
$ws.Range("K110").Value = 3877.7334
$ws.Range("L110").Value = 4178.6
$ws.Range("M110").Value = -1832.7334
$ws.Range("N110").Value = -8268.6
$ws.Range("H122").Value = 2372.8948
$ws.Range("I122").Value = 2376.9443
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 7130.8329
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -4680.8329
$ws.Range("N122").Value = -11800
$ws.Range("H136").Value = 3604.7778
$ws.Range("J136").Value = 4382.0557
$ws.Range("L136").Value = 13146.1671
$ws.Range("N136").Value = -18246.1671

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1349.6
$ws.Range("J64").Value = 1329.3334
$ws.Range("L64").Value = 1329.3334
$ws.Range("N64").Value = -1779.3334
$ws.Range("H67").Value = 1349.6
$ws.Range("J67").Value = 1329.3334
$ws.Range("L67").Value = 1329.3334
$ws.Range("N67").Value = -2889.3334
$ws.Range("H86").Value = 2506.8
$ws.Range("I86").Value = 2341
$ws.Range("J86").Value = 3999
$ws.Range("K86").Value = 2341
$ws.Range("L86").Value = 3999
$ws.Range("M86").Value = -1218
$ws.Range("N86").Value = -6245
$ws.Range("H89").Value = 2506.8
$ws.Range("I89").Value = 2341
$ws.Range("J89").Value = 3999
$ws.Range("K89").Value = 11705
$ws.Range("L89").Value = 19995
$ws.Range("M89").Value = -6089
$ws.Range("N89").Value = -31227
$ws.Range("H110").Value = 90000
$ws.Range("J110").Value = 90000
$ws.Range("L110").Value = 90000
$ws.Range("N110").Value = -98180
$ws.Range("H134").Value = 3781.5
$ws.Range("I134").Value = 3619.75
$ws.Range("J134").Value = 3902.8125
$ws.Range("K134").Value = 10859.25
$ws.Range("L134").Value = 11708.4375
$ws.Range("M134").Value = -8324.25
$ws.Range("N134").Value = -16778.4375

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5296.4
$ws.Range("I16").Value = 5296.4
$ws.Range("K16").Value = 5296.4
$ws.Range("M16").Value = -5009.4
$ws.Range("H31").Value = 8036.6284
$ws.Range("I31").Value = 3999.1428
$ws.Range("J31").Value = 10728.286
$ws.Range("K31").Value = 3999.1428
$ws.Range("L31").Value = 10728.286
$ws.Range("M31").Value = -3704.1428
$ws.Range("N31").Value = -11318.286
$ws.Range("H34").Value = 8036.6284
$ws.Range("I34").Value = 3999.1428
$ws.Range("J34").Value = 10728.286
$ws.Range("K34").Value = 3999.1428
$ws.Range("L34").Value = 10728.286
$ws.Range("M34").Value = -3797.1428
$ws.Range("N34").Value = -11132.286
$ws.Range("H62").Value = 10999
$ws.Range("I62").Value = 12236.75
$ws.Range("J62").Value = 6048
$ws.Range("K62").Value = 12236.75
$ws.Range("L62").Value = 6048
$ws.Range("M62").Value = -11612.75
$ws.Range("N62").Value = -7296
$ws.Range("H65").Value = 10999
$ws.Range("I65").Value = 12236.75
$ws.Range("J65").Value = 6048
$ws.Range("K65").Value = 61183.75
$ws.Range("L65").Value = 30240
$ws.Range("M65").Value = -58063.75
$ws.Range("N65").Value = -36480
$ws.Range("H113").Value = 5296.4
$ws.Range("I113").Value = 5296.4
$ws.Range("K113").Value = 5296.4
$ws.Range("M113").Value = -3126.4

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 73.333336
$ws.Range("I2").Value = 97.5
$ws.Range("K2").Value = 585
$ws.Range("M2").Value = -472
$ws.Range("H9").Value = 79999.75
$ws.Range("J9").Value = 79999.75
$ws.Range("L9").Value = 239999.25
$ws.Range("N9").Value = -240447.25

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1489.3334
$ws.Range("I97").Value = 1519.8667
$ws.Range("J97").Value = 1336.6666
$ws.Range("K97").Value = 1519.8667
$ws.Range("L97").Value = 1336.6666
$ws.Range("M97").Value = -1023.8667
$ws.Range("N97").Value = -2328.6666
$ws.Range("H98").Value = 32025.666
$ws.Range("J98").Value = 32025.666
$ws.Range("L98").Value = 32025.666
$ws.Range("N98").Value = -38015.666
$ws.Range("H113").Value = 2699.9375
$ws.Range("I113").Value = 2823
$ws.Range("J113").Value = 2166.6667
$ws.Range("K113").Value = 2823
$ws.Range("L113").Value = 2166.6667
$ws.Range("M113").Value = -653
$ws.Range("N113").Value = -6506.6667
$ws.Range("H132").Value = 4454.4546
$ws.Range("I132").Value = 4610.4
$ws.Range("K132").Value = 13831.2
$ws.Range("M132").Value = -11301.2

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 887.7222
$ws.Range("I22").Value = 826.63635
$ws.Range("J22").Value = 983.7143
$ws.Range("K22").Value = 826.63635
$ws.Range("L22").Value = 983.7143
$ws.Range("M22").Value = -531.63635
$ws.Range("N22").Value = -1573.7143
$ws.Range("H27").Value = 887.7222
$ws.Range("I27").Value = 826.63635
$ws.Range("J27").Value = 983.7143
$ws.Range("K27").Value = 826.63635
$ws.Range("L27").Value = 983.7143
$ws.Range("M27").Value = -719.63635
$ws.Range("N27").Value = -1197.7143
$ws.Range("H40").Value = 5612.1665
$ws.Range("I40").Value = 6154.6
$ws.Range("J40").Value = 2900
$ws.Range("K40").Value = 6154.6
$ws.Range("L40").Value = 2900
$ws.Range("M40").Value = -6018.6
$ws.Range("N40").Value = -3172
$ws.Range("H46").Value = 1930.24
$ws.Range("I46").Value = 1000.36365
$ws.Range("K46").Value = 1000.36365
$ws.Range("M46").Value = -812.36365
$ws.Range("H61").Value = 4350.75
$ws.Range("I61").Value = 3485.6155
$ws.Range("K61").Value = 3485.6155
$ws.Range("M61").Value = -3283.6155
$ws.Range("H80").Value = 33888
$ws.Range("I80").Value = 32940.234
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 32940.234
$ws.Range("L80").Value = 50000
$ws.Range("M80").Value = -31817.234
$ws.Range("N80").Value = -52246
$ws.Range("H83").Value = 33888
$ws.Range("I83").Value = 32940.234
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 98820.70199999999
$ws.Range("L83").Value = 150000
$ws.Range("M83").Value = -93204.70199999999
$ws.Range("N83").Value = -161232
$ws.Range("H93").Value = 4062.25
$ws.Range("I93").Value = 4062.25
$ws.Range("K93").Value = 4062.25
$ws.Range("M93").Value = -2814.25
$ws.Range("H113").Value = 4350.75
$ws.Range("I113").Value = 3485.6155
$ws.Range("K113").Value = 3485.6155
$ws.Range("M113").Value = -1315.6155

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10803.182
$ws.Range("I41").Value = 12179.8
$ws.Range("K41").Value = 12179.8
$ws.Range("M41").Value = -11789.8
$ws.Range("H81").Value = 2520.111
$ws.Range("I81").Value = 2520.111
$ws.Range("K81").Value = 5040.222
$ws.Range("M81").Value = -3979.222
$ws.Range("H84").Value = 2520.111
$ws.Range("I84").Value = 2520.111
$ws.Range("K84").Value = 25201.11
$ws.Range("M84").Value = -19897.11
$ws.Range("H126").Value = 3629.2
$ws.Range("I126").Value = 3832.4443
$ws.Range("K126").Value = 11497.3329
$ws.Range("M126").Value = -9027.332900000001
$ws.Range("H132").Value = 2798.4285
$ws.Range("I132").Value = 2431.8333
$ws.Range("K132").Value = 7295.499899999999
$ws.Range("M132").Value = -4765.499899999999
$ws.Range("H136").Value = 1974.3572
$ws.Range("I136").Value = 2058.8235
$ws.Range("J136").Value = 1843.8182
$ws.Range("K136").Value = 6176.470499999999
$ws.Range("L136").Value = 5531.4546
$ws.Range("M136").Value = -3626.470499999999
$ws.Range("N136").Value = -10631.4546
